# Applies the "Forgot few things to add" commit:
#  - splits a few single-run bullet items into multiple runs (as the
#    author's Word session apparently did while editing them)
#  - fixes stray spaces / a typo ("shour" -> "hour")
#  - appends a batch of new November bullet items that were missing

$d = $word.ActiveDocument
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

function Get-ParaIndexByText($doc, [string]$text) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        if ($doc.Paragraphs($i).Range.Text.Trim() -eq $text) {
            return $i
        }
    }
    return -1
}

function Replace-ParaXml($doc, [string]$oldText, [string]$innerXml) {
    $idx = Get-ParaIndexByText $doc $oldText
    if ($idx -eq -1) {
        throw "Paragraph not found: $oldText"
    }
    $rng = $doc.Paragraphs($idx).Range
    $rng.InsertXML($innerXml)
}

# ---------------------------------------------------------------------
# 1) "Participated in all class meeting so far in October"
# ---------------------------------------------------------------------
$xml1 = @"
<w:p $wNs>
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:numPr><w:ilvl w:val="0"/><w:numId w:val="4"/></w:numPr>
    <w:rPr><w:lang w:val="en-US"/></w:rPr>
  </w:pPr>
  <w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">Participated in all stand up meetings </w:t></w:r>
  <w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>during</w:t></w:r>
  <w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> October</w:t></w:r>
</w:p>
"@
Replace-ParaXml $d "Participated in all class meeting so far in October" $xml1

# ---------------------------------------------------------------------
# 2) "Created test plan document ( October 10/19/2020 2 hours)"
# ---------------------------------------------------------------------
$xml2 = @"
<w:p $wNs>
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:numPr><w:ilvl w:val="0"/><w:numId w:val="4"/></w:numPr>
    <w:rPr><w:lang w:val="en-US"/></w:rPr>
  </w:pPr>
  <w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">Created test plan document </w:t></w:r>
  <w:proofErr w:type="gramStart"/>
  <w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>( October</w:t></w:r>
  <w:proofErr w:type="gramEnd"/>
  <w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> 10/19/2020 2 hours)</w:t></w:r>
</w:p>
"@
Replace-ParaXml $d "Created test plan document ( October 10/19/2020 2 hours)" $xml2

# ---------------------------------------------------------------------
# 3) "Found paper with nose key points which will be used in project ( 11/2/2020 2 hours)"
# ---------------------------------------------------------------------
$xml3 = @"
<w:p $wNs>
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:numPr><w:ilvl w:val="0"/><w:numId w:val="4"/></w:numPr>
    <w:rPr><w:lang w:val="en-US"/></w:rPr>
  </w:pPr>
  <w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Found paper with nose key points which will be used in project (</w:t></w:r>
  <w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>1</w:t></w:r>
  <w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>1/2/2020 2 hours)</w:t></w:r>
</w:p>
"@
Replace-ParaXml $d "Found paper with nose key points which will be used in project ( 11/2/2020 2 hours)" $xml3

# ---------------------------------------------------------------------
# 4) "Found better shape predictor to extract more nose points ( 11/3/2020 2.5 hours)"
# ---------------------------------------------------------------------
$xml4 = @"
<w:p $wNs>
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:numPr><w:ilvl w:val="0"/><w:numId w:val="4"/></w:numPr>
    <w:rPr><w:lang w:val="en-US"/></w:rPr>
  </w:pPr>
  <w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Found better shape predictor to extract more nose points (11/3/2020 2.5 hours)</w:t></w:r>
</w:p>
"@
Replace-ParaXml $d "Found better shape predictor to extract more nose points ( 11/3/2020 2.5 hours)" $xml4

# ---------------------------------------------------------------------
# 5) "Created new points on nose( 11/6/2020 2 hours)"
# ---------------------------------------------------------------------
$xml5 = @"
<w:p $wNs>
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:numPr><w:ilvl w:val="0"/><w:numId w:val="4"/></w:numPr>
    <w:rPr><w:lang w:val="en-US"/></w:rPr>
  </w:pPr>
  <w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Created new points on nose</w:t></w:r>
  <w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>
  <w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>(11/6/2020 2 hours)</w:t></w:r>
</w:p>
"@
Replace-ParaXml $d "Created new points on nose( 11/6/2020 2 hours)" $xml5

# ---------------------------------------------------------------------
# 6) "Extended test cases and worked on test plan ( 11/7/2020 2 hours)"
# ---------------------------------------------------------------------
$xml6 = @"
<w:p $wNs>
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:numPr><w:ilvl w:val="0"/><w:numId w:val="4"/></w:numPr>
    <w:rPr><w:lang w:val="en-US"/></w:rPr>
  </w:pPr>
  <w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Extended test cases and worked on test plan (11/7/2020 2 hours)</w:t></w:r>
</w:p>
"@
Replace-ParaXml $d "Extended test cases and worked on test plan ( 11/7/2020 2 hours)" $xml6

# ---------------------------------------------------------------------
# 7) "Updated traceability matrix( 11/9/2020 1 shour)"
# ---------------------------------------------------------------------
$xml7 = @"
<w:p $wNs>
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:numPr><w:ilvl w:val="0"/><w:numId w:val="4"/></w:numPr>
    <w:rPr><w:lang w:val="en-US"/></w:rPr>
  </w:pPr>
  <w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Updated</w:t></w:r>
  <w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> traceability matrix</w:t></w:r>
  <w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> (</w:t></w:r>
  <w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>11/</w:t></w:r>
  <w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>9</w:t></w:r>
  <w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>/2020 1 hour)</w:t></w:r>
</w:p>
"@
Replace-ParaXml $d "Updated traceability matrix( 11/9/2020 1 shour)" $xml7

# ---------------------------------------------------------------------
# 8) "Created new database for the app ( 11/9/2020 3 hours)" + 8 new bullets
# ---------------------------------------------------------------------
$xml8 = @"
<w:p $wNs>
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:numPr><w:ilvl w:val="0"/><w:numId w:val="4"/></w:numPr>
    <w:rPr><w:lang w:val="en-US"/></w:rPr>
  </w:pPr>
  <w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:lastRenderedPageBreak/><w:t>Created new database for the app (11/9/2020 3 hours)</w:t></w:r>
</w:p>
<w:p $wNs>
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:numPr><w:ilvl w:val="0"/><w:numId w:val="4"/></w:numPr>
    <w:rPr><w:lang w:val="en-US"/></w:rPr>
  </w:pPr>
  <w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Found new shape predictor for better face calculations (11/13/2020 1 hour)</w:t></w:r>
</w:p>
<w:p $wNs>
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:numPr><w:ilvl w:val="0"/><w:numId w:val="4"/></w:numPr>
    <w:rPr><w:lang w:val="en-US"/></w:rPr>
  </w:pPr>
  <w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Worked on SRS document, reviewed comments, worked on models (11/15/2020 3 hours)</w:t></w:r>
</w:p>
<w:p $wNs>
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:numPr><w:ilvl w:val="0"/><w:numId w:val="4"/></w:numPr>
    <w:rPr><w:lang w:val="en-US"/></w:rPr>
  </w:pPr>
  <w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Worked on SDS document,</w:t></w:r>
  <w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>
  <w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>reviewed comments, worked on models (11/15/2020 3 hours)</w:t></w:r>
</w:p>
<w:p $wNs>
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:numPr><w:ilvl w:val="0"/><w:numId w:val="4"/></w:numPr>
    <w:rPr><w:lang w:val="en-US"/></w:rPr>
  </w:pPr>
  <w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Updated SRS document, Updated models, requirements, and descriptions (11/16/2020 4 hours)</w:t></w:r>
</w:p>
<w:p $wNs>
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:numPr><w:ilvl w:val="0"/><w:numId w:val="4"/></w:numPr>
    <w:rPr><w:lang w:val="en-US"/></w:rPr>
  </w:pPr>
  <w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">Updated </w:t></w:r>
  <w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>SDD</w:t></w:r>
  <w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> document, updated models, added interface design and descriptions (11/16/2020 4 hours)</w:t></w:r>
</w:p>
<w:p $wNs>
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:numPr><w:ilvl w:val="0"/><w:numId w:val="4"/></w:numPr>
    <w:rPr><w:lang w:val="en-US"/></w:rPr>
  </w:pPr>
  <w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Found formula for scaling pixels to mm (11/18/2020 2 hours of research)</w:t></w:r>
</w:p>
<w:p $wNs>
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:numPr><w:ilvl w:val="0"/><w:numId w:val="4"/></w:numPr>
    <w:rPr><w:lang w:val="en-US"/></w:rPr>
  </w:pPr>
  <w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Updated models for SRS and SDD (11/28/2020 2 hours)</w:t></w:r>
</w:p>
<w:p $wNs>
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:numPr><w:ilvl w:val="0"/><w:numId w:val="4"/></w:numPr>
    <w:rPr><w:lang w:val="en-US"/></w:rPr>
  </w:pPr>
  <w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Participated in all stand up meetings during November</w:t></w:r>
</w:p>
"@
Replace-ParaXml $d "Created new database for the app ( 11/9/2020 3 hours)" $xml8

Write-Output "Edits applied."
